# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" colour scheme (only used by the Notes Master)
#   ppt/theme/theme2.xml -> "Integral" colour scheme (used by the Slide Master / the
#                            presentation's visible design)
#
# The authored change swaps the two themes' content: the design that slides actually
# use should now carry the "Office Theme" palette (dk1/lt1/dk2/lt2/accent1-6/hlink/
# folHlink), while what used to be the "Integral" palette moves to the other theme
# part. The font scheme and format scheme are identical between the two themes
# already, so only the twelve theme colours need to change.
#
# PowerPoint's object model exposes the live theme colours for the deck's design via
# Slide.ThemeColorScheme (SlideMaster/Design "Integral"); rewrite every slot to the
# "Office Theme" RGB values so the on-disk theme part matches the target palette.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$colors = $s.ThemeColorScheme

$colors.Item(1).RGB  = 0         # dk1      -> 000000
$colors.Item(2).RGB  = 16777215  # lt1      -> FFFFFF
$colors.Item(3).RGB  = 6968388   # dk2      -> 44546A
$colors.Item(4).RGB  = 15132391  # lt2      -> E7E6E6
$colors.Item(5).RGB  = 13998939  # accent1  -> 5B9BD5
$colors.Item(6).RGB  = 3243501   # accent2  -> ED7D31
$colors.Item(7).RGB  = 10855845  # accent3  -> A5A5A5
$colors.Item(8).RGB  = 49407     # accent4  -> FFC000
$colors.Item(9).RGB  = 12874308  # accent5  -> 4472C4
$colors.Item(10).RGB = 4697456   # accent6  -> 70AD47
$colors.Item(11).RGB = 12673797  # hlink    -> 0563C1
$colors.Item(12).RGB = 7491477   # folHlink -> 954F72
